# Apply the changes described by the diff:
# 1. Bump Version from 1.0.3 to 1.0.4 on the "isa_template" sheet.
# 2. Rename "Parameter [chromatography instrument]" -> "Component [chromatography instrument]"
#    in both the worksheet header row and the table column definition.
# 3. Rename "Parameter [chromatography column model]" -> "Component [chromatography column model]"
#    in both the worksheet header row and the table column definition.
# 4. Rename "Output [Raw Data File]" -> "Output [Data]" in both the worksheet header row
#    and the table column definition.
# 5. Update the NCIT term accession URL from the purl.obolibrary.org form to the bioregistry.io form.

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("isa_template")
$tableSheet = $wb.Worksheets.Item("New Table")

# 1. Version bump
$metaSheet.Range("B4").Value = "1.0.4"

# 2-4. Header row renames on the table sheet
$tableSheet.Range("K1").Value = "Component [chromatography instrument]"
$tableSheet.Range("N1").Value = "Component [chromatography column model]"
$tableSheet.Range("W1").Value = "Output [Data]"

# 5. Update NCIT accession number URL value in the data row
$tableSheet.Range("D2").Value = "https://bioregistry.io/NCIT:C30014"

# Keep the ListObject (table) column names in sync with the header row cells.
$table = $tableSheet.ListObjects.Item(1)
$table.ListColumns.Item("Component [chromatography instrument]").Name = "Component [chromatography instrument]"
$table.ListColumns.Item("Component [chromatography column model]").Name = "Component [chromatography column model]"
$table.ListColumns.Item("Output [Data]").Name = "Output [Data]"

$wb.Save()
